$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# The XPath-style strings for the child form references gain a leading
# slash (absolute paths), matching the updated linking logic.
$ws.Range("F3").Value = "/child/a/name"
$ws.Range("F4").Value = "/child/a/age"
$ws.Range("F6").Value = "/child/extra_info"

# Move the active selection from F8 to F7 on the survey sheet.
[void]$ws.Range("F7").Select()
